# Update examples and code revision
# - Swap the 'type' (B) and 'description' (E) columns on the Processes sheet
#   (cgam_processes now spans A1:D1 instead of A1:E1).
# - Reduce the waste recycle ratio for QG on WasteDefinition from 0.9 to 0.5
#   and leave WasteDefinition as the active/selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Defined name 'cgam_processes' now only covers columns A:D (was A:E)
# ---------------------------------------------------------------------
$wb.Names("cgam_processes").RefersTo = "=Processes!`$A`$1:`$D`$1"

# ---------------------------------------------------------------------
# 2. Processes sheet: swap column B (type) and column E (description)
# ---------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("Processes")

for ($r = 1; $r -le 11; $r++) {
    $bVal = $wsProc.Cells.Item($r, 2).Value()
    $eVal = $wsProc.Cells.Item($r, 5).Value()
    $wsProc.Cells.Item($r, 2).Value = $eVal
    $wsProc.Cells.Item($r, 5).Value = $bVal
}

# Column E is a bit wider now that it holds descriptions
$wsProc.Columns.Item(5).ColumnWidth = 14.75

# Selection moves onto the (now) description column
[void]$wsProc.Activate()
[void]$wsProc.Range("E1:E11").Select()

# ---------------------------------------------------------------------
# 3. WasteDefinition sheet: QG recycle ratio 0.9 -> 0.5
# ---------------------------------------------------------------------
$wsWaste = $wb.Worksheets.Item("WasteDefinition")
$wsWaste.Range("C3").Value = 0.5

# WasteDefinition ends up as the active sheet/selection
[void]$wsWaste.Activate()
[void]$wsWaste.Range("C4").Select()
